# Updated cryptos list (price/volume refresh + a few rank swaps), matching
# the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.947.41'
$ws.Range("E2").Value = '  -5.53%  '
$ws.Range("D3").Value = '3.367.97'
$ws.Range("E3").Value = '  -7.30%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''183.98'
$ws.Range("E5").Value = '  -9.86%  '
$ws.Range("D6").Value = '''525.79'
$ws.Range("E6").Value = '  -8.34%  '
$ws.Range("D7").Value = '''0.598'
$ws.Range("E7").Value = '  -4.62%  '
$ws.Range("D8").Value = '3.360.17'
$ws.Range("E8").Value = '  -7.30%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '''0.620'
$ws.Range("E10").Value = '  -10.13%  '
$ws.Range("D11").Value = '''57.11'
$ws.Range("E11").Value = '  -7.92%  '
$ws.Range("E12").Value = '  -13.78%  '
$ws.Range("D13").Value = '''0.0000252'
$ws.Range("E13").Value = '  -12.95%  '
$ws.Range("D14").Value = '''9.19'
$ws.Range("E14").Value = '  -9.84%  '
$ws.Range("D15").Value = '3.912.77'
$ws.Range("E15").Value = '  -7.12%  '
$ws.Range("E16").Value = '  -3.51%  '
$ws.Range("D17").Value = '3.371.05'
$ws.Range("E17").Value = '  -7.24%  '
$ws.Range("D18").Value = '64.641.03'
$ws.Range("E18").Value = '  -5.69%  '
$ws.Range("D19").Value = '''17.28'
$ws.Range("E19").Value = '  -9.47%  '
$ws.Range("D20").Value = '''11.02'
$ws.Range("E20").Value = '  -11.98%  '
$ws.Range("D21").Value = '''0.961'
$ws.Range("E21").Value = '  -11.45%  '
$ws.Range("D22").Value = '''369.97'
$ws.Range("E22").Value = '  -9.61%  '
$ws.Range("D23").Value = '''80.69'
$ws.Range("E23").Value = '  -6.41%  '
$ws.Range("D24").Value = '''3.69'
$ws.Range("E24").Value = '  -12.49%  '
$ws.Range("D25").Value = '''10.79'
$ws.Range("E25").Value = '  -16.91%  '
$ws.Range("D26").Value = '''3.75'
$ws.Range("E26").Value = '  -5.79%  '
$ws.Range("D27").Value = '''5.82'
$ws.Range("E27").Value = '  -5.63%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '''2.62'
$ws.Range("E28").Value = '  -11.60%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''11.36'
$ws.Range("E29").Value = '  -10.77%  '
$ws.Range("D30").Value = '''8.45'
$ws.Range("E30").Value = '  -10.58%  '
$ws.Range("D31").Value = '''29.36'
$ws.Range("E31").Value = '  -7.76%  '
$ws.Range("D32").Value = '''660.26'
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("E33").Value = '  -17.43%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '''61.05'
$ws.Range("E34").Value = '  -4.57%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").Value = '''11.11'
$ws.Range("E35").Value = '  -10.46%  '
$ws.Range("D36").Value = '''0.104'
$ws.Range("E36").Value = '  -10.14%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = '''36.35'
$ws.Range("E38").Value = '  -14.31%  '
$ws.Range("D39").Value = '''0.377'
$ws.Range("E39").Value = '  -11.07%  '
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''0.126'
$ws.Range("E41").Value = '  -7.76%  '
$ws.Range("D42").Value = '2.826.34'
$ws.Range("E42").Value = '  -12.00%  '
$ws.Range("D43").Value = '''2.73'
$ws.Range("E43").Value = '  -16.65%  '
$ws.Range("D44").Value = '0.0₃0627'
$ws.Range("E44").Value = '  -22.70%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '''2.62'
$ws.Range("E45").Value = '  -9.44%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0389'
$ws.Range("E46").Value = '  -7.91%  '
$ws.Range("D47").Value = '''2.31'
$ws.Range("E47").Value = '  -15.81%  '
$ws.Range("D48").Value = '''137.21'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("E49").Value = '  -6.65%  '
$ws.Range("D50").Value = '''2.87'
$ws.Range("E50").Value = '  -7.25%  '
$ws.Range("D51").Value = '''2.56'
$ws.Range("E51").Value = '  -7.23%  '
